$d = $word.ActiveDocument
$replacements = @(
    @("2025-01-30 Thursday", "2025-01-31 Friday"),
    @("47+3=50", "14+47=61"),
    @("88-36=52", "93-5=88"),
    @("41+1=42", "67-19=48"),
    @("88-46=42", "11-0=11"),
    @("81+18=99", "75-21=54"),
    @("80+17=97", "74-60=14"),
    @("36+4=40", "95-89=6"),
    @("57+40=97", "97-24=73"),
    @("48-44=4", "34-28=6"),
    @("33-13=20", "16+81=97"),
    @("3+82=85", "1+47=48"),
    @("67+9=76", "75-42=33"),
    @("27-2=25", "33-5=28"),
    @("84-5=79", "60-35=25"),
    @("36+14=50", "53+28=81"),
    @("87-43=44", "24-11=13"),
    @("49+41=90", "31+0=31"),
    @("75+18=93", "29-14=15"),
    @("48+46=94", "44-44=0"),
    @("89-82=7", "67+10=77"),
    @("34+18=52", "79+5=84"),
    @("94-88=6", "11+33=44"),
    @("95-72=23", "14+4=18"),
    @("59-41=18", "33+42=75"),
    @("46+12=58", "96-29=67"),
    @("1+0=1", "3+94=97"),
    @("76-7=69", "71+15=86"),
    @("41-4=37", "51+10=61"),
    @("9+11=20", "13+71=84"),
    @("12+10=22", "61+23=84"),
    @("24+66=90", "75+8=83"),
    @("73-29=44", "1+33=34"),
    @("70+10=80", "87-73=14"),
    @("51+17=68", "87-59=28"),
    @("44+29=73", "7-0=7"),
    @("33-32=1", "46-22=24"),
    @("48-12=36", "48+35=83"),
    @("35+61=96", "97-32=65"),
    @("5+17=22", "19+24=43"),
    @("90-18=72", "97-63=34"),
    @("81-77=4", "9+59=68"),
    @("77-68=9", "37+30=67"),
    @("5-5=0", "76-73=3"),
    @("54-48=6", "42-13=29"),
    @("92-15=77", "11+80=91"),
    @("87+7=94", "18+5=23"),
    @("49+31=80", "9+40=49"),
    @("46+6=52", "24+20=44"),
    @("99-6=93", "47+33=80"),
    @("99-7=92", "14+15=29"),
    @("75-5=70", "21+21=42"),
    @("19+63=82", "75-56=19"),
    @("98-11=87", "98-1=97"),
    @("37+14=51", "92-18=74"),
    @("86-66=20", "77-45=32"),
    @("31-5=26", "4+83=87"),
    @("15+34=49", "16+9=25"),
    @("71+17=88", "63+9=72"),
    @("40+24=64", "83-74=9"),
    @("71-22=49", "54+20=74"),
    @("53+6=59", "91-89=2"),
    @("6+26=32", "95-43=52"),
    @("21+59=80", "18-14=4"),
    @("37+36=73", "93-24=69"),
    @("41+49=90", "99-17=82"),
    @("3+52=55", "40-24=16"),
    @("12+55=67", "23+72=95"),
    @("15-9=6", "36+9=45"),
    @("23+29=52", "74+14=88"),
    @("68-6=62", "86-27=59"),
    @("1+97=98", "11+18=29"),
    @("3+34=37", "37+17=54"),
    @("35+18=53", "98-27=71"),
    @("48-24=24", "43+23=66"),
    @("24+59=83", "15+21=36"),
    @("95-52=43", "47-9=38"),
    @("48-31=17", "28+56=84"),
    @("70-25=45", "77-52=25"),
    @("89+9=98", "81+0=81"),
    @("13+84=97", "27+65=92"),
    @("48+19=67", "92-40=52"),
    @("91-11=80", "68-51=17"),
    @("14+55=69", "56+0=56"),
    @("72-48=24", "19+67=86"),
    @("33-4=29", "98-6=92"),
    @("41+32=73", "89-64=25"),
    @("18+47=65", "44+26=70"),
    @("58-30=28", "83-62=21"),
    @("53-37=16", "58-56=2"),
    @("80+18=98", "97-14=83"),
    @("94-71=23", "33+18=51"),
    @("22+57=79", "11+17=28"),
    @("66-10=56", "67-9=58"),
    @("54+42=96", "24+10=34"),
    @("50+5=55", "35+63=98"),
    @("65+15=80", "62+11=73"),
    @("81+5=86", "46+11=57"),
    @("65+29=94", "69+8=77"),
    @("5+0=5", "85-70=15"),
    @("93-71=22", "0+91=91"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replaced $($replacements.Count) text runs"
